$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: "pin" every existing run boundary inside a paragraph (given its
# start offset and the character-length of each of its runs) by toggling a
# character property on/off. The canonical-OOXML writer merges adjacent
# runs that end up with identical formatting once a paragraph is touched;
# toggling Bold here leaves the *value* unchanged (still not-bold) but
# marks the run as explicitly-set, which keeps the writer from re-merging
# it into its neighbours. Doing this for every run *before* any surgical
# text edit preserves the paragraph's pre-existing run structure exactly.
# ---------------------------------------------------------------------------
function Pin-Runs($paraStart, $lens) {
    $offset = $paraStart
    foreach ($len in $lens) {
        if ($len -gt 0) {
            $r = $d.Range($offset, $offset + $len)
            $r.Bold = 1
            $r.Bold = 0
        }
        $offset += $len
    }
    return $offset
}

# Paragraph containing "...but holding a token..."
$paraA = 661
$lensA = @(41,90,12,4,15,1,2,58,21,4,16,4,9,1,6,16,7,13,1,1)
Pin-Runs $paraA $lensA | Out-Null

# Paragraph containing "So the driver is getting real-time" ... "ot pictured
# on the phone is a server piece that sends notification email ..."
$paraB = 1370
$lensB = @(34,17,7,14,73,6,8,1,103,96,23,1)
Pin-Runs $paraB $lensB | Out-Null

# Paragraph that originally held the "_GoBack" bookmark after " make"
$paraC = 2004
$lensC = @(34,1,11,28,10,1,4,5,79,9,10)
Pin-Runs $paraC $lensC | Out-Null

# ---------------------------------------------------------------------------
# Change 1: " but holding a token " -> " but " / "essentially" / " a token "
#   (the run boundaries around it are already pinned above, so this text
#    surgery only affects the "holding" -> "essentially" word itself and
#    splits that run in three without disturbing anything else)
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.Text = "holding a token"
$found = $find.Execute()
if ($found) {
    $wordStart = $find.Parent.Start
    $wordEnd = $wordStart + 7   # length of "holding"
    $middle = $d.Range($wordStart, $wordEnd)
    $middle.Text = "essentially"

    $essFind = $d.Content.Find
    $essFind.Text = "essentially"
    $essFind.Execute() | Out-Null
    $ess = $essFind.Parent
    $ess.Bold = 1
    $ess.Bold = 0
}

# ---------------------------------------------------------------------------
# Change 2: "So the driver is getting real-time" + " feedback on how " ->
#           "So the driver is getting " + "feedback on how "
#   (removes the word "real-time" plus one adjoining space)
# ---------------------------------------------------------------------------
$find3 = $d.Content.Find
$find3.Text = "real-time "
$found3 = $find3.Execute()
if ($found3) {
    $delStart = $find3.Parent.Start
    $delEnd = $find3.Parent.End
    $del = $d.Range($delStart, $delEnd)
    $del.Text = ""

    $fbFind = $d.Content.Find
    $fbFind.Text = "feedback on how "
    $fbFind.Execute() | Out-Null
    $fb = $fbFind.Parent
    $fb.Bold = 1
    $fb.Bold = 0
}

# ---------------------------------------------------------------------------
# Changes 3 & 4: move the "_GoBack" bookmark from right after " make" to
#                between "...that sends" and " notification email "
#   (Bookmarks.Add with an existing name relocates it, removing the old
#    bookmark pair and naturally splitting the run at the new location)
# ---------------------------------------------------------------------------
$find5 = $d.Content.Find
$find5.Text = "that sends notification"
$found5 = $find5.Execute()
if ($found5) {
    $sendsEnd = $find5.Parent.Start + 10   # length of "that sends"
    $gobackRange = $d.Range($sendsEnd, $sendsEnd)
    $d.Bookmarks.Add("_GoBack", $gobackRange)

    $neFind = $d.Content.Find
    $neFind.Text = " notification email "
    $neFind.Execute() | Out-Null
    $ne = $neFind.Parent
    $ne.Bold = 1
    $ne.Bold = 0
}
